$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Scroll the view so B4 is the top-left visible cell (best effort - matches
# the xlsx diff's sheetView topLeftCell="B4")
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 2

# --- Shared string order matters: insert new strings in this exact order so the
# resulting sharedStrings table matches: Qrr(60), Switching Charge (Coul.)(61),
# FET Driver(62), Source/Sink Current(63), Turn on time (S)(64), MIC4102(65),
# Design for MIC4102...(66)

# 1) J19: "Qrr" (was blank styled cell s="1" -> becomes plain text "Qrr")
$ws.Range("J19").Value = "Qrr"

# 2) J17: text changes from "Qoss(Coul.)" to "Switching Charge (Coul.)"
$ws.Range("J17").Value = "Switching Charge (Coul.)"
# K17 formula: 299*POWER(10,-9) -> 26*POWER(10,-9)
$ws.Range("K17").Formula = "=26*POWER(10,-9)"

# 3) New J18/K18 cells re-introduce the old "Qoss(Coul.)" label + old 299e-9 formula
$ws.Range("J18").Value = "Qoss(Coul.)"
$ws.Range("K18").Formula = "=299*POWER(10,-9)"

# 4) K16 formula: 100*POWER(10,-9) -> 80*POWER(10,-9)
$ws.Range("K16").Formula = "=80*POWER(10,-9)"

# 5) New M7 (merged M7:P7) "FET Driver" label, styled like K7/L7 (centered, s=7)
$ws.Range("K7").Copy()
$ws.Range("M7:P7").PasteSpecial(-4122)
$ws.Range("M7").Value = "FET Driver"
$ws.Range("M7:P7").Merge()

# 6) New M11 "Source/Sink Current" label and N11 value (3)
$ws.Range("M11").Value = "Source/Sink Current"
$ws.Range("N11").Value = 3

# 7) New J20 "Turn on time (S)" label and K20 formula =K17/N11
$ws.Range("J20").Value = "Turn on time (S)"
$ws.Range("K20").Formula = "=K17/N11"

# 8) New M8 "MIC4102" and O8 note text
$ws.Range("M8").Value = "MIC4102"
$ws.Range("O8").Value = "Design for MIC4102 for PWM input, add a not gate that can be DNP normally, but populated if need to change dot MIC4103"

# Autofit column M to approximate the bestFit width Excel applied to the new column
$ws.Columns("M:M").AutoFit()

# K26 formula: B18*E13*10*POWER(10,-9)*E8 -> B18*E13*K20*E8
$ws.Range("K26").Formula = "=B18*E13*K20*E8"

# K28 formula: (K17/2)*B18*E8 -> (K18/2)*B18*E8
$ws.Range("K28").Formula = "=(K18/2)*B18*E8"

# --- Comments ---------------------------------------------------------
$c1 = $ws.Range("J16").AddComment("Shelby R:" + [char]10 + "Excellent resource on gate charge characteristics: " + [char]10 + "https://www.microsemi.com/document-portal/doc_view/14697-making-use-of-gate-charge-information-in-mosfet-and-igbt-data-sheets")

$c2 = $ws.Range("J17").AddComment("Shelby R:" + [char]10 + "This is the charge that is relevant for switching losses = Qgd + (Qg-Qgd-Qgth)=Qsw")

$c3 = $ws.Range("J19").AddComment("Helpful page on benefits of low Qrr:" + [char]10 + "https://efficiencywins.nexperia.com/efficient-products/qrr-overlooked-and-underappreciated-in-efficiency-battle.html")

# Final selection matches the diff's sheetView selection (activeCell/sqref = M11)
$ws.Range("M11").Select()
